$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 4833.5
$ws.Range("I17").Value = 3000.5
$ws.Range("J17").Value = 5750
$ws.Range("K17").Value = 9001.5
$ws.Range("L17").Value = 17250
$ws.Range("M17").Value = -8833.5
$ws.Range("N17").Value = -17586

$ws.Range("H40").Value = 800
$ws.Range("I40").Value = 800
$ws.Range("K40").Value = 800
$ws.Range("M40").Value = -625

$ws.Range("H58").Value = 974.5
$ws.Range("I58").Value = 637.25
$ws.Range("J58").Value = 1649
$ws.Range("K58").Value = 1911.75
$ws.Range("L58").Value = 4947
$ws.Range("M58").Value = -1761.75
$ws.Range("N58").Value = -5247

$ws.Range("H100").Value = 4158.2
$ws.Range("J100").Value = 6329.6665
$ws.Range("L100").Value = 6329.6665
$ws.Range("N100").Value = -7411.6665

$ws.Range("H103").Value = 500
$ws.Range("J103").Value = 500
$ws.Range("L103").Value = 1500
$ws.Range("N103").Value = -2672

$ws.Range("H106").Value = 3599.75
$ws.Range("I106").Value = 3599.75
$ws.Range("K106").Value = 3599.75
$ws.Range("M106").Value = -2968.75

$ws.Range("H132").Value = 2961.75
$ws.Range("I132").Value = 3313.8572
$ws.Range("J132").Value = 497
$ws.Range("K132").Value = 9941.571599999999
$ws.Range("L132").Value = 1491
$ws.Range("M132").Value = -7411.571599999999
$ws.Range("N132").Value = -6551

$ws.Range("H137").Value = 2237.75
$ws.Range("I137").Value = 1843.1428
$ws.Range("K137").Value = 5529.428400000001
$ws.Range("M137").Value = -2979.428400000001

$ws.Range("H138").Value = 1465.381
$ws.Range("I138").Value = 597.875
$ws.Range("K138").Value = 1793.625
$ws.Range("M138").Value = 3346.375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2652.889
$ws.Range("I122").Value = 2652.889
$ws.Range("K122").Value = 7958.667
$ws.Range("M122").Value = -5508.667

$ws.Range("H132").Value = 2996.5
$ws.Range("I132").Value = 2597.8
$ws.Range("J132").Value = 4990
$ws.Range("K132").Value = 7793.400000000001
$ws.Range("L132").Value = 14970
$ws.Range("M132").Value = -5263.400000000001
$ws.Range("N132").Value = -20030

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2525.238
$ws.Range("I86").Value = 1997.5834
$ws.Range("J86").Value = 3228.7778
$ws.Range("K86").Value = 1997.5834
$ws.Range("L86").Value = 3228.7778
$ws.Range("M86").Value = -874.5834
$ws.Range("N86").Value = -5474.7778

$ws.Range("H89").Value = 2525.238
$ws.Range("I89").Value = 1997.5834
$ws.Range("J89").Value = 3228.7778
$ws.Range("K89").Value = 9987.916999999999
$ws.Range("L89").Value = 16143.889
$ws.Range("M89").Value = -4371.916999999999
$ws.Range("N89").Value = -27375.889

$ws.Range("H122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()

$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 23714.834
$ws.Range("I58").Value = 15948.5
$ws.Range("K58").Value = 15948.5
$ws.Range("M58").Value = -15745.5

$ws.Range("H62").Value = 1898
$ws.Range("I62").Value = 1898
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 1898
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -1274
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 1898
$ws.Range("I65").Value = 1898
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 9490
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -6370
$ws.Range("N65").ClearContents()

$ws.Range("H92").Value = 47666.668
$ws.Range("J92").Value = 47666.668
$ws.Range("L92").Value = 47666.668
$ws.Range("N92").Value = -52658.668

$ws.Range("H134").Value = 6967.5557
$ws.Range("I134").Value = 7979
$ws.Range("J134").Value = 5703.25
$ws.Range("K134").Value = 23937
$ws.Range("L134").Value = 17109.75
$ws.Range("M134").Value = -21402
$ws.Range("N134").Value = -22179.75

$ws.Range("H136").Value = 23714.834
$ws.Range("I136").Value = 15948.5
$ws.Range("K136").Value = 47845.5
$ws.Range("M136").Value = -45295.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H29").Value = 51.304348
$ws.Range("J29").Value = 40.909092
$ws.Range("L29").Value = 122.727276
$ws.Range("N29").Value = -676.727276

$ws.Range("H34").Value = 1000
$ws.Range("I34").Value = 1000
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 3000
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -2916
$ws.Range("N34").ClearContents()

$ws.Range("H36").Value = 648.5
$ws.Range("I36").Value = 397
$ws.Range("J36").Value = 900
$ws.Range("K36").Value = 1191
$ws.Range("L36").Value = 2700
$ws.Range("M36").Value = -1022
$ws.Range("N36").Value = -3038

$ws.Range("H39").Value = 2250
$ws.Range("J39").Value = 4000
$ws.Range("L39").Value = 12000
$ws.Range("N39").Value = -12588

$ws.Range("H46").Value = 1500
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 1500
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 4500
$ws.Range("N46").Value = -4682
$ws.Range("M46").ClearContents()

$ws.Range("H50").Value = 474
$ws.Range("I50").Value = 531.4286
$ws.Range("K50").Value = 1594.2858
$ws.Range("M50").Value = -1113.2858

$ws.Range("H53").Value = 474
$ws.Range("I53").Value = 531.4286
$ws.Range("K53").Value = 1594.2858
$ws.Range("M53").Value = -1113.2858

$ws.Range("H55").Value = 0
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("M55").ClearContents()
$ws.Range("N55").ClearContents()

$ws.Range("H60").Value = 30
$ws.Range("I60").Value = 30
$ws.Range("K60").Value = 90
$ws.Range("M60").Value = 161

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 6968050.5
$ws.Range("I122").Value = 9647017
$ws.Range("J122").Value = 2738.6
$ws.Range("K122").Value = 28941051
$ws.Range("L122").Value = 8215.799999999999
$ws.Range("M122").Value = -28938601
$ws.Range("N122").Value = -13115.8

$ws.Range("H126").Value = 2359.3333
$ws.Range("I126").Value = 2039
$ws.Range("J126").Value = 3000
$ws.Range("K126").Value = 6117
$ws.Range("L126").Value = 9000
$ws.Range("M126").Value = -3647
$ws.Range("N126").Value = -13940

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3524
$ws.Range("I122").Value = 3145
$ws.Range("J122").Value = 3966.1667
$ws.Range("K122").Value = 9435
$ws.Range("L122").Value = 11898.5001
$ws.Range("M122").Value = -6985
$ws.Range("N122").Value = -16798.5001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 16365.429
$ws.Range("I4").Value = 35001.668
$ws.Range("J4").Value = 2388.25
$ws.Range("K4").Value = 35001.668
$ws.Range("L4").Value = 2388.25
$ws.Range("M4").Value = -34888.668
$ws.Range("N4").Value = -2614.25

$ws.Range("H81").Value = 978.4
$ws.Range("I81").Value = 972.5
$ws.Range("K81").Value = 1945
$ws.Range("M81").Value = -884

$ws.Range("H84").Value = 978.4
$ws.Range("I84").Value = 972.5
$ws.Range("K84").Value = 9725
$ws.Range("M84").Value = -4421

$ws.Range("H100").Value = 777.4
$ws.Range("I100").Value = 747
$ws.Range("K100").Value = 1494
$ws.Range("M100").Value = -953

$ws.Range("H107").Value = 515.4
$ws.Range("I107").Value = 515.4
$ws.Range("K107").Value = 1546.2
$ws.Range("M107").Value = 373.8000000000002
